$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grammar fix: "Identifie" -> "Identify" across the sprint log entries.
$ws.Range("D4").Value = "Identify 3 code smells"
$ws.Range("D5").Value = "Identify 3 code smells"
$ws.Range("D6").Value = "Identify 3 code smells"
$ws.Range("D7").Value = "Identify 3 code smells"
$ws.Range("D8").Value = "Identify 3 code smells"

$ws.Range("D10").Value = "Identify 3 design paterns"
$ws.Range("D11").Value = "Identify 3 design paterns"
$ws.Range("D12").Value = "Identify 3 design paterns"
$ws.Range("D13").Value = "Identify 3 design paterns"
$ws.Range("D14").Value = "Identify 3 design paterns "

# Update the active cell selection to match the saved view state.
$ws.Range("D14").Select()
